# Fruta / hortaliza, semanal
# Insert a new weekly record at row 962 on the single data sheet, pushing the
# existing rows 962-1006 down to 963-1007.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 962, shifting everything below
# (through the former last row 1006) down by one.
$ws.Rows.Item(962).Insert()

# Populate the newly inserted row 962 with the new weekly price record.
$ws.Cells.Item(962, 1).Value = 10
$ws.Cells.Item(962, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(962, 3).Value = "La Araucanía"
$ws.Cells.Item(962, 4).Value = 44939
$ws.Cells.Item(962, 5).Value = 9
$ws.Cells.Item(962, 6).Value = "Fruta"
$ws.Cells.Item(962, 7).Value = 100104
$ws.Cells.Item(962, 8).Value = "Frutos de pepita"
$ws.Cells.Item(962, 9).Value = 100104005
$ws.Cells.Item(962, 10).Value = "Pera"
$ws.Cells.Item(962, 11).Value = "Favorita De Clapp"
$ws.Cells.Item(962, 12).Value = "Primera"
$ws.Cells.Item(962, 13).Value = 80
$ws.Cells.Item(962, 14).Value = 20000
$ws.Cells.Item(962, 15).Value = 20000
$ws.Cells.Item(962, 16).Value = 20000
$ws.Cells.Item(962, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(962, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(962, 19).Value = 1111
$ws.Cells.Item(962, 20).Value = 18
